$wb = $excel.ActiveWorkbook

# Update the "想去人数" (interest count) column F for each sheet
# per the diff. Column F is the 6th column.

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 1781
$ws.Cells.Item(3, 6).Value = 10409
$ws.Cells.Item(5, 6).Value = 17
$ws.Cells.Item(6, 6).Value = 610
$ws.Cells.Item(7, 6).Value = 85
$ws.Cells.Item(8, 6).Value = 1714
$ws.Cells.Item(9, 6).Value = 436
$ws.Cells.Item(10, 6).Value = 7
$ws.Cells.Item(11, 6).Value = 245
$ws.Cells.Item(13, 6).Value = 534
$ws.Cells.Item(15, 6).Value = 149
$ws.Cells.Item(17, 6).Value = 1027
$ws.Cells.Item(18, 6).Value = 34
$ws.Cells.Item(19, 6).Value = 120
$ws.Cells.Item(20, 6).Value = 412
$ws.Cells.Item(21, 6).Value = 412
$ws.Cells.Item(23, 6).Value = 354
$ws.Cells.Item(24, 6).Value = 51
$ws.Cells.Item(25, 6).Value = 1073
$ws.Cells.Item(26, 6).Value = 1121
$ws.Cells.Item(27, 6).Value = 1215
$ws.Cells.Item(28, 6).Value = 215
$ws.Cells.Item(29, 6).Value = 1418
$ws.Cells.Item(30, 6).Value = 722
$ws.Cells.Item(31, 6).Value = 258
$ws.Cells.Item(32, 6).Value = 30
$ws.Cells.Item(34, 6).Value = 661
$ws.Cells.Item(35, 6).Value = 257
$ws.Cells.Item(36, 6).Value = 743
$ws.Cells.Item(38, 6).Value = 781
$ws.Cells.Item(39, 6).Value = 805
$ws.Cells.Item(42, 6).Value = 834
$ws.Cells.Item(44, 6).Value = 1377
$ws.Cells.Item(45, 6).Value = 52
$ws.Cells.Item(46, 6).Value = 729
$ws.Cells.Item(48, 6).Value = 720
$ws.Cells.Item(49, 6).Value = 85

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(19, 6).Value = 1148
$ws.Cells.Item(21, 6).Value = 2258
$ws.Cells.Item(22, 6).Value = 1122
$ws.Cells.Item(23, 6).Value = 350
$ws.Cells.Item(24, 6).Value = 696
$ws.Cells.Item(25, 6).Value = 97
$ws.Cells.Item(27, 6).Value = 30
$ws.Cells.Item(30, 6).Value = 377
$ws.Cells.Item(35, 6).Value = 166
$ws.Cells.Item(36, 6).Value = 198
$ws.Cells.Item(41, 6).Value = 138
$ws.Cells.Item(43, 6).Value = 13
$ws.Cells.Item(46, 6).Value = 86

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 829
$ws.Cells.Item(5, 6).Value = 214
$ws.Cells.Item(6, 6).Value = 2575
$ws.Cells.Item(7, 6).Value = 4258
$ws.Cells.Item(10, 6).Value = 412
$ws.Cells.Item(11, 6).Value = 363
$ws.Cells.Item(12, 6).Value = 267
$ws.Cells.Item(13, 6).Value = 159
$ws.Cells.Item(14, 6).Value = 78

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 1781
$ws.Cells.Item(3, 6).Value = 829
$ws.Cells.Item(5, 6).Value = 10409
$ws.Cells.Item(6, 6).Value = 214
$ws.Cells.Item(7, 6).Value = 4258
$ws.Cells.Item(9, 6).Value = 363
$ws.Cells.Item(10, 6).Value = 1714
$ws.Cells.Item(11, 6).Value = 436
$ws.Cells.Item(12, 6).Value = 245
$ws.Cells.Item(15, 6).Value = 149
$ws.Cells.Item(17, 6).Value = 1027
$ws.Cells.Item(18, 6).Value = 34
$ws.Cells.Item(19, 6).Value = 412
$ws.Cells.Item(20, 6).Value = 412
$ws.Cells.Item(22, 6).Value = 354
$ws.Cells.Item(23, 6).Value = 51
$ws.Cells.Item(24, 6).Value = 2258
$ws.Cells.Item(25, 6).Value = 2258
$ws.Cells.Item(26, 6).Value = 1122
$ws.Cells.Item(27, 6).Value = 1073
$ws.Cells.Item(28, 6).Value = 1121
$ws.Cells.Item(29, 6).Value = 1215
$ws.Cells.Item(30, 6).Value = 97
$ws.Cells.Item(31, 6).Value = 1418
$ws.Cells.Item(32, 6).Value = 722
$ws.Cells.Item(33, 6).Value = 377
$ws.Cells.Item(34, 6).Value = 661
$ws.Cells.Item(36, 6).Value = 743
$ws.Cells.Item(38, 6).Value = 781
$ws.Cells.Item(40, 6).Value = 805
$ws.Cells.Item(42, 6).Value = 834
$ws.Cells.Item(44, 6).Value = 1377
$ws.Cells.Item(45, 6).Value = 52
$ws.Cells.Item(46, 6).Value = 138
$ws.Cells.Item(48, 6).Value = 729
$ws.Cells.Item(49, 6).Value = 720
$ws.Cells.Item(50, 6).Value = 86
